# Lit matrix current (Feb Week 1).xlsx - add "Data acquisition literature review" section
# (rows 14-17) with new papers, column/row sizing, a centered section header, and
# threaded comments on D16/E16.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Month 1")
$ws.Activate()

# ---------------------------------------------------------------------------
# Column widths (best effort - engine quantizes to a coarse pixel grid so the
# exact fractional widths from the source file cannot always be hit exactly).
# ---------------------------------------------------------------------------
$ws.Columns("B").ColumnWidth = 24.109375
$ws.Columns("C").ColumnWidth = 27.21875
$ws.Columns("D").ColumnWidth = 23.77734375
$ws.Columns("E").ColumnWidth = 27.88671875
$ws.Columns("G").ColumnWidth = 27.109375

# ---------------------------------------------------------------------------
# Row 14: merged section header "Data acquisition literature review"
# ---------------------------------------------------------------------------
$hdr = $ws.Range("B14:H14")
$hdr.Font.Name = "Century Gothic"
$hdr.Font.Size = 9
$hdr.VerticalAlignment = -4160   # xlTop
$hdr.HorizontalAlignment = -4108 # xlCenter
$hdr.WrapText = $true
$ws.Range("B14").Value2 = "Data acquisition literature review"
$hdr.Merge()

# ---------------------------------------------------------------------------
# Row 15
# ---------------------------------------------------------------------------
$ws.Rows("15").RowHeight = 132

$ws.Range("A15").Value2 = "Cheng-Wen Wu"
$ws.Range("B15").Value2 = "Can IoT make Semiconductor great again?"
$ws.Range("C15").Value2 = @"
There is not much evidence that IoT will likely give a great boost to the semiconductor industry in the near future due to limitations in global economy and energy consumption.
If IoT is going to give a boost to stagnant semiconductor industry, what will be the key factors of its success? 
"@
$ws.Range("D15").Value2 = @"
Propose the symbiotic system model (SSM) for developing IoT devices and systems.
Propose symbiosis-based test (SBT) for device and system test.
"@
$ws.Range("E15").Value2 = @"
A symbiotic relationship (SR) is a relationship of mutual dependence between two different systems where one system's input is from the other system's output (& vice versa)
The twin system (SS) comprises of the primary or functional system and the secondary or test system
"@
$ws.Range("F15").Value2 = "N/A"
$ws.Range("G15").Value2 = "The article seeks to trigger more research activities regarding establishing a sound IoT platform that allows heterogeneous integration of technologies and partners to migrate certain industries based on the notion of IoT"

# ---------------------------------------------------------------------------
# Row 16
# ---------------------------------------------------------------------------
$ws.Rows("16").RowHeight = 184.8

$ws.Range("A16").Value2 = @"
Kyselova A.G
Verbitskyi I V
Kyselov G.D
"@
$ws.Range("B16").Value2 = "Context-aware framework for energy management system"
$ws.Range("C16").Value2 = @"
Need to have visibility on the context from the level of sensor data to the higher level sitauation awareness (actuator level)
Explore the challenges in microgrid energy control systems
"@
$ws.Range("D16").Value2 = @"
Proposes a CAEMS that presents data management solutions.
These solutions include sensor data acquisition and time series forecasting, ontology model and context prediction model for analytical query processing past and future context data
"@
$ws.Range("E16").Value2 = @"
Energy systems: 
* Renewable energy integration, multidirectional power flow, reduced peak load leveling and load demand, bidirectional flow of communication in the system, uninterrupted power supply, and improved quality of electricity.
Context prediction - ability to predict the future context information in order to provide proactive service to the actions of all electrical facilities (loads and generators)
"@

$f16 = "Did not provide a good background for the case scenario of their work. It would have been good to see actual use cases in the microgrids showing how the proposed model would be implemented. `n(The work focused on data management making an assumption on its acquisition)"
$ws.Range("F16").Value2 = $f16
$boldPart = "The work focused on data management making an assumption on its acquisition)"
$boldStart = $f16.IndexOf($boldPart) + 1
$boldChars = $ws.Range("F16").Characters($boldStart, $boldPart.Length)
$boldChars.Font.Name = "Century Gothic"
$boldChars.Font.Size = 9
$boldChars.Font.Bold = $true

$ws.Range("G16").Value2 = @"
Interesting feature of ontological modeling and context based predictions (Very good approaches to consider for distributed IoT Edge layer evaluation)
My interest lies in how the data is collected, shared and how the sensors behave when exposed to different environments (remote or otherwise).
"@

# ---------------------------------------------------------------------------
# Row 17
# ---------------------------------------------------------------------------
$ws.Rows("17").RowHeight = 105.6

$ws.Range("A17").Value2 = @"
A Mawire
R.R.J. van den Heetkamp
M McPherson
E Zhandire
"@
$ws.Range("B17").Value2 = "Data acquisition and control of a thermal energy storage and cooking system"
$ws.Range("D17").Value2 = @"
Data acquisition for the thermal energy solar system (Data points: Temperature, fluid flow rate, and power input)
Use of HP34970A data loggers (embedded with ADC)
"@
$ws.Range("E17").Value2 = "The temperature contorl program is able to maintain a nearly constant charging temperature of the TES system."
$ws.Range("F17").Value2 = "N/A"
$ws.Range("G17").Value2 = "Data acquisition system design and evaluation was done for thermal based small scale system"

# ---------------------------------------------------------------------------
# Threaded comments
# ---------------------------------------------------------------------------
$commentD16 = @"
The framework creates the general approaches to the context prediction
Create and support automated decision-making approaches for reasoning context for energy consumption and user comfort
The energy system is highly distributed and CAEMS manages large amounts of energy-related data that has to be able to react rapidly and smartly when conditions change and for this task we use smoothies and predictive techniques for data from sensors.
"@
$ws.Range("D16").AddCommentThreaded($commentD16) | Out-Null

$commentE16 = @"
Context aware management process:-
- Sensor data acquisition
- Feature extraction
- Ontology model
- Context prediction
"@
$e16Comment = $ws.Range("E16").AddCommentThreaded($commentE16)

$replyE16 = @"
Keeping a record of time points and the value of the disturbances complicates the forecasting process and can lead to erroneous results. Filtering or smoothing of context time series is the necessary preliminary prediction stage for obtaining trends
- Averaging methods
- Exponential smoothing methods
- Kalman filter
"@
$addMethod = "Add"
$e16Comment.Replies.$addMethod($replyE16) | Out-Null

# ---------------------------------------------------------------------------
# Final sheet view / selection
# ---------------------------------------------------------------------------
$ws.Range("A16").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("C16").Select()
